$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update participation/eligibility-driven metrics for rows 2-6 (Year 1-5)
# seeded from recalculated participation flags based on existing deferral rates

# Row 2
$ws.Range("C2").Value = 102
$ws.Range("D2").Value = 85
$ws.Range("E2").Value = 0.8333333333333334
$ws.Range("F2").Value = 0.8333333333333334
$ws.Range("G2").Value = 0.1003684358597651
$ws.Range("H2").Value = 0.08364036321647093
$ws.Range("I2").Value = 456350.2764925769
$ws.Range("J2").Value = 165561.1384472884
$ws.Range("L2").Value = 165561.1384472884
$ws.Range("M2").Value = 621911.4149398654
$ws.Range("N2").Value = 10078372.3488
$ws.Range("O2").Value = 9670631.418699998
$ws.Range("P2").Value = 0.01642736869778395
$ws.Range("Q2").Value = 0.01711999261259659

# Row 3
$ws.Range("D3").Value = 86
$ws.Range("E3").Value = 0.8349514563106796
$ws.Range("F3").Value = 0.8349514563106796
$ws.Range("G3").Value = 0.09909943690069946
$ws.Range("H3").Value = 0.08274321915980733
$ws.Range("I3").Value = 475621.2913790991
$ws.Range("J3").Value = 172598.2258965795
$ws.Range("L3").Value = 172598.2258965795
$ws.Range("M3").Value = 648219.5172756788
$ws.Range("N3").Value = 10579921.128564
$ws.Range("O3").Value = 10172547.970561
$ws.Range("P3").Value = 0.01631375355252824
$ws.Range("Q3").Value = 0.01696705942268081

# Row 4
$ws.Range("C4").Value = 104
$ws.Range("D4").Value = 87
$ws.Range("E4").Value = 0.8365384615384616
$ws.Range("F4").Value = 0.8365384615384616
$ws.Range("G4").Value = 0.09886985879837833
$ws.Range("H4").Value = 0.08270843957172033
$ws.Range("I4").Value = 503855.2163199777
$ws.Range("J4").Value = 179628.8157732298
$ws.Range("L4").Value = 179628.8157732298
$ws.Range("M4").Value = 683484.0320932076
$ws.Range("N4").Value = 10893476.50522092
$ws.Range("O4").Value = 10485132.15247783
$ws.Range("P4").Value = 0.01648957664590721
$ws.Range("Q4").Value = 0.01713176459399991

# Row 5
$ws.Range("D5").Value = 91
$ws.Range("E5").Value = 0.8666666666666667
$ws.Range("F5").Value = 0.8666666666666667
$ws.Range("G5").Value = 0.09652624012210713
$ws.Range("H5").Value = 0.08365607477249284
$ws.Range("I5").Value = 536050.1427417491
$ws.Range("J5").Value = 193152.3152125126
$ws.Range("L5").Value = 193152.3152125126
$ws.Range("M5").Value = 729202.4579542616
$ws.Range("N5").Value = 11344132.16627755
$ws.Range("O5").Value = 10933437.48295217
$ws.Range("P5").Value = 0.01702662772095447
$ws.Range("Q5").Value = 0.01766620200771103

# Row 6
$ws.Range("D6").Value = 90
$ws.Range("E6").Value = 0.8490566037735849
$ws.Range("F6").Value = 0.8490566037735849
$ws.Range("G6").Value = 0.09703484391241361
$ws.Range("H6").Value = 0.08238807501997383
$ws.Range("I6").Value = 551966.9472694583
$ws.Range("J6").Value = 198049.5347916165
$ws.Range("L6").Value = 198049.5347916165
$ws.Range("M6").Value = 750016.4820610748
$ws.Range("N6").Value = 11800478.14216588
$ws.Range("O6").Value = 11386012.61834073
$ws.Range("P6").Value = 0.01678317881747004
$ws.Range("Q6").Value = 0.01739410814217752
